$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ISO-Alpha-3 code column (B) and the Country Name column (C) were
# entered in the wrong order when the table was first built - swap them
# back so B holds the ISO Alpha-3 code and C holds the full country name,
# matching the column headers already in row 1.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Clear the stray multi-column selection left over from editing and put
# the cursor back on A1, as a finished/clean workbook would have.
[void]$ws.Range("A1").Select()
